$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 -> becomes the data that was previously in row 3 (NT / Spillkråka / Dryocopus martius)
$ws.Range("A2").Value = 131054136
$ws.Range("B2").Value = 57881
$ws.Range("D2").Value = "NT"
$ws.Range("E2").Value = 100049
$ws.Range("F2").Value = "Spillkråka"
$ws.Range("G2").Value = "Dryocopus martius"
$ws.Range("H2").Value = "(Linnaeus, 1758)"
$ws.Range("Q2").Value = 660622
$ws.Range("R2").Value = 6661065
$ws.Range("Z2").Value = "12:23"
$ws.Range("AB2").Value = "12:23"
$ws.Range("AC2").Value = ""

# Row 3 -> becomes the data that was previously in row 2 (LC / Tjäder / Tetrao urogallus)
$ws.Range("A3").Value = 131054138
$ws.Range("B3").Value = 57073
$ws.Range("D3").Value = "LC"
$ws.Range("E3").Value = 100138
$ws.Range("F3").Value = "Tjäder"
$ws.Range("G3").Value = "Tetrao urogallus"
$ws.Range("H3").Value = "Linnaeus, 1758"
$ws.Range("Q3").Value = 660731
$ws.Range("R3").Value = 6661229
$ws.Range("Z3").Value = "12:04"
$ws.Range("AB3").Value = "12:04"
$ws.Range("AC3").Value = "2 tuppar"
